$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 32329.666
$ws.Range("J7").Value = 32329.666
$ws.Range("L7").Value = 32329.666
$ws.Range("N7").Value = -32553.666
$ws.Range("H14").Value = 32329.666
$ws.Range("J14").Value = 32329.666
$ws.Range("L14").Value = 32329.666
$ws.Range("N14").Value = -32711.666
$ws.Range("H31").Value = 4766.7144
$ws.Range("I31").Value = 3728
$ws.Range("J31").Value = 10999
$ws.Range("K31").Value = 11184
$ws.Range("L31").Value = 32997
$ws.Range("M31").Value = -10954
$ws.Range("N31").Value = -33457
$ws.Range("H33").Value = 740.38464
$ws.Range("I33").Value = 808.3333
$ws.Range("K33").Value = 808.3333
$ws.Range("M33").Value = -579.3333
$ws.Range("H40").Value = 5833.3335
$ws.Range("H116").Value = 9753.764999999999
$ws.Range("I116").Value = 10113
$ws.Range("J116").Value = 4006
$ws.Range("K116").Value = 10113
$ws.Range("L116").Value = 4006
$ws.Range("M116").Value = -6671
$ws.Range("N116").Value = -10890
$ws.Range("H137").Value = 2016.9445
$ws.Range("J137").Value = 2500.889
$ws.Range("L137").Value = 7502.667
$ws.Range("N137").Value = -12602.667

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 3666.3333
$ws.Range("J17").Value = 3666.3333
$ws.Range("L17").Value = 3666.3333
$ws.Range("N17").Value = -4012.3333
$ws.Range("H32").Value = 7884.0977
$ws.Range("I32").Value = 5986.385
$ws.Range("K32").Value = 5986.385
$ws.Range("M32").Value = -5699.385
$ws.Range("H122").Value = 2760.4546
$ws.Range("I122").Value = 3188.5715
$ws.Range("J122").Value = 2011.25
$ws.Range("K122").Value = 9565.7145
$ws.Range("L122").Value = 6033.75
$ws.Range("M122").Value = -7115.7145
$ws.Range("N122").Value = -10933.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 3306.6
$ws.Range("I29").Value = 508
$ws.Range("K29").Value = 508
$ws.Range("M29").Value = -219
$ws.Range("H33").Value = 1250.5
$ws.Range("I33").Value = 1
$ws.Range("J33").Value = 2500
$ws.Range("K33").Value = 1
$ws.Range("L33").Value = 2500
$ws.Range("M33").Value = 335
$ws.Range("N33").Value = -3172

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 20000
$ws.Range("J25").Value = 25000
$ws.Range("L25").Value = 25000
$ws.Range("N25").Value = -25348
$ws.Range("H35").Value = 2967.7334
$ws.Range("I35").Value = 649.5
$ws.Range("J35").Value = 7604.2
$ws.Range("K35").Value = 649.5
$ws.Range("L35").Value = 7604.2
$ws.Range("M35").Value = -355.5
$ws.Range("N35").Value = -8192.200000000001
$ws.Range("H132").Value = 125001250
$ws.Range("I132").Value = 142857860
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 428573580
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -428571050
$ws.Range("N132").Value = -20060

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 216.9
$ws.Range("I17").Value = 161.5
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 484.5
$ws.Range("L17").Value = 900
$ws.Range("M17").Value = -315.5
$ws.Range("N17").Value = -1238
$ws.Range("H32").Value = 3000
$ws.Range("J32").Value = 3000
$ws.Range("L32").Value = 9000
$ws.Range("N32").Value = -9566
$ws.Range("H46").Value = 186
$ws.Range("I46").Value = 202
$ws.Range("J46").Value = 154
$ws.Range("K46").Value = 606
$ws.Range("L46").Value = 462
$ws.Range("M46").Value = -515
$ws.Range("N46").Value = -644
$ws.Range("H68").Value = 1182.3334
$ws.Range("J68").Value = 1999
$ws.Range("L68").Value = 5997
$ws.Range("N68").Value = -7619
$ws.Range("H69").Value = 782.75
$ws.Range("I69").Value = 716
$ws.Range("K69").Value = 2148
$ws.Range("M69").Value = -1337
$ws.Range("H71").Value = 1182.3334
$ws.Range("J71").Value = 1999
$ws.Range("L71").Value = 17991
$ws.Range("N71").Value = -26103
$ws.Range("H72").Value = 782.75
$ws.Range("I72").Value = 716
$ws.Range("K72").Value = 6444
$ws.Range("M72").Value = -2388
$ws.Range("H109").Value = 3992.7856
$ws.Range("I109").Value = 1474.75
$ws.Range("J109").Value = 5000
$ws.Range("K109").Value = 4424.25
$ws.Range("L109").Value = 15000
$ws.Range("M109").Value = -3384.25
$ws.Range("N109").Value = -17080
$ws.Range("H113").Value = 35248.207
$ws.Range("I113").Value = 111552.445
$ws.Range("K113").Value = 334657.335
$ws.Range("M113").Value = -332487.335
$ws.Range("H121").Value = 87841.07000000001
$ws.Range("I121").Value = 204199.6
$ws.Range("J121").Value = 23197.445
$ws.Range("K121").Value = 612598.8
$ws.Range("L121").Value = 69592.33499999999
$ws.Range("M121").Value = -611288.8
$ws.Range("N121").Value = -72212.33499999999
$ws.Range("H131").Value = 2057.5
$ws.Range("H132").Value = 1994.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 8171.143
$ws.Range("I22").Value = 875
$ws.Range("J22").Value = 17899.334
$ws.Range("K22").Value = 875
$ws.Range("L22").Value = 17899.334
$ws.Range("M22").Value = -346
$ws.Range("N22").Value = -18957.334
$ws.Range("H46").Value = 11500
$ws.Range("J46").Value = 18000
$ws.Range("L46").Value = 18000
$ws.Range("N46").Value = -18312
$ws.Range("H122").Value = 206831.5
$ws.Range("I122").Value = 402363.34
$ws.Range("J122").Value = 11299.667
$ws.Range("K122").Value = 1207090.02
$ws.Range("L122").Value = 33899.001
$ws.Range("M122").Value = -1204640.02
$ws.Range("N122").Value = -38799.001
$ws.Range("H126").Value = 7407.3335
$ws.Range("I126").Value = 7796.375
$ws.Range("K126").Value = 23389.125
$ws.Range("M126").Value = -20919.125
$ws.Range("H132").Value = 5954236.5
$ws.Range("I132").Value = 7354367.5
$ws.Range("J132").Value = 3678.5
$ws.Range("K132").Value = 22063102.5
$ws.Range("L132").Value = 11035.5
$ws.Range("M132").Value = -22060572.5
$ws.Range("N132").Value = -16095.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 16000
$ws.Range("I39").Value = 16000
$ws.Range("K39").Value = 16000
$ws.Range("M39").Value = -15540
$ws.Range("H132").Value = 12001845
$ws.Range("I132").Value = 15001420
$ws.Range("J132").Value = 3545.375
$ws.Range("K132").Value = 45004260
$ws.Range("L132").Value = 10636.125
$ws.Range("M132").Value = -45001730
$ws.Range("N132").Value = -15696.125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 22194.25
$ws.Range("I29").Value = 33333
$ws.Range("J29").Value = 11055.5
$ws.Range("K29").Value = 33333
$ws.Range("L29").Value = 11055.5
$ws.Range("M29").Value = -33043
$ws.Range("N29").Value = -11635.5
$ws.Range("H132").Value = 11117036
$ws.Range("I132").Value = 14289970
$ws.Range("J132").Value = 11766.7
$ws.Range("K132").Value = 42869910
$ws.Range("L132").Value = 35300.10000000001
$ws.Range("M132").Value = -42867380
$ws.Range("N132").Value = -40360.10000000001
